$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.941706120967865
$ws.Range("B1").Value = 1.807476162910461
$ws.Range("C1").Value = 4.5207839012146
$ws.Range("D1").Value = 2.926390409469604
$ws.Range("E1").Value = 0.4090909957885742
